# Updated cryptos list on Thu May 16 11:10:52 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "66.146.78"
$ws.Range("E2").Value = "  +5.40%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.995.25"
$ws.Range("E3").Value = "  +3.04%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB (numeric-looking value, use text prefix to keep it textual)
$ws.Range("D5").Value = "'580.77"

# Row 6 - Solana
$ws.Range("D6").Value = "'163.13"
$ws.Range("E6").Value = "  +12.80%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.04%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +3.14%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.991.62"

# Row 10 - Toncoin
$ws.Range("D10").Value = "'6.51"
$ws.Range("E10").Value = "  -5.65%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +2.70%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "'0.455"
$ws.Range("E12").Value = "  +5.05%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  +5.38%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'34.53"
$ws.Range("E14").Value = "  +4.94%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  -0.88%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "66.153.86"
$ws.Range("E16").Value = "  +5.54%  "

# Row 17 - WrappedliquidstakedEther2.0
$ws.Range("D17").Value = "3.491.98"
$ws.Range("E17").Value = "  +3.02%  "

# Row 18 - Polkadot
$ws.Range("E18").Value = "  +4.00%  "

# Row 19 - WrappedEther
$ws.Range("D19").Value = "2.998.84"
$ws.Range("E19").Value = "  +3.24%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'452.25"
$ws.Range("E20").Value = "  +5.23%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "'13.81"
$ws.Range("E21").Value = "  +5.35%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  +3.82%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  +6.23%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'82.18"
$ws.Range("E24").Value = "  +4.20%  "

# Row 25 - Fetch.AI
$ws.Range("E25").Value = "  +13.74%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").Value = "'12.26"
$ws.Range("E26").Value = "  +2.59%  "

# Row 27 - RenderToken
$ws.Range("D27").Value = "'10.38"

# Row 28 - Dai
$ws.Range("E28").Value = "  +0.02%  "

# Row 29 - NEARProtocol
$ws.Range("D29").Value = "'8.10"
$ws.Range("E29").Value = "  +13.15%  "

# Row 30 - ImmutableX
$ws.Range("E30").Value = "  +18.46%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +4.86%  "

# Row 32 - PEPE
$ws.Range("D32").Value = "'0.0000103"

# Row 33 - EthereumClassic
$ws.Range("D33").Value = "'27.22"
$ws.Range("E33").Value = "  +5.22%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  +3.43%  "

# Row 35 - FirstDigitalUSD
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.08%  "

# Row 36 - Mantle
$ws.Range("D36").Value = "'0.991"
$ws.Range("E36").Value = "  +3.86%  "

# Row 37 - Filecoin
$ws.Range("E37").Value = "  +7.44%  "

# Row 38 - Stacks
$ws.Range("D38").Value = "'2.07"
$ws.Range("E38").Value = "  +8.37%  "

# Row 39 - OKB
$ws.Range("D39").Value = "'49.52"
$ws.Range("E39").Value = "  +1.65%  "

# Row 40 - dogwifhat
$ws.Range("E40").Value = "  -0.29%  "

# Row 41 - TheGraph
$ws.Range("D41").Value = "'0.309"
$ws.Range("E41").Value = "  +15.40%  "

# Row 42 & 43 - Kaspa and Arweave swap places with updated values
$ws.Range("B42").Value = "Arweave"
$ws.Range("C42").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D42").Value = "'44.02"
$ws.Range("E42").Value = "  +6.99%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.122"
$ws.Range("E43").Value = "  +6.79%  "

# Row 44 - Cosmos
$ws.Range("D44").Value = "'8.44"
$ws.Range("E44").Value = "  +4.26%  "

# Row 45 - Bittensor
$ws.Range("D45").Value = "'400.31"
$ws.Range("E45").Value = "  +11.48%  "

# Row 46 - VeChain
$ws.Range("E46").Value = "  +5.41%  "

# Row 47 - Maker
$ws.Range("D47").Value = "2.778.01"
$ws.Range("E47").Value = "  +2.34%  "

# Row 48 - Monero
$ws.Range("D48").Value = "'133.34"
$ws.Range("E48").Value = "  -0.01%  "

# Row 50 - InjectiveProtocol
$ws.Range("D50").Value = "'23.80"
$ws.Range("E50").Value = "  +10.90%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  +3.34%  "
